# Update "Orders" sheet: add Number (F21) to the existing last row, then append
# new order-line rows 22-31 (PackageID in column A, FlowerName in column C,
# Number in column F). Numeric-looking identifiers/quantities are written as
# text, consistent with the rest of the sheet (columns A and F).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "7"

$ws.Range("C22").Value = "175_火灵鸟_Free Spirit_Rosa rugosa Thunb._20stems"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "7"

$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "8"
$ws.Range("C23").Value = "197_粉红雪山_Sweet Avalanche_Rosa rugosa Thunb._20stems"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "6"

$ws.Range("C24").Value = "152_白荔枝_White Ohara_Rosa rugosa Thunb._20stems"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = "6"

$ws.Range("C25").Value = "274_仙子之吻_undefined_Rosa rugosa Thunb._10stems"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "32"

$ws.Range("C26").Value = "302_彩星 浅粉_Tinted Gypso light pink_undefined_0.5kg"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "40"

$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "9"
$ws.Range("C27").Value = "110_绣球单瓣浅蓝_Hydrangea Light Blue S_Hydrangea L._1stem"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "40"

$ws.Range("C28").Value = "106_绣球单瓣粉_Hydrangea Pink S_Hydrangea L._1stem"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "40"

$ws.Range("C29").Value = "100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem"
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "115"

$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = "10"
$ws.Range("C30").Value = "157_流沙_Quicksand_Rosa rugosa Thunb._20stems"
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "60"

$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = "11"

# Update "Summary" sheet: TotalNumber (G2) reflects the newly appended orders
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Value = "01519618194232023211522524141030811776632404040115600"
